$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVTR")

# Row 4: Inventory
$ws.Range("B4").Value = 778000000.0
$ws.Range("C4").Value = 740000000.0
$ws.Range("D4").Value = 730000000.0
$ws.Range("E4").Value = 738000000.0
$ws.Range("F4").Value = 686000000.0

# Row 13: Accounts Payable
$ws.Range("B13").Value = 706000000.0
$ws.Range("C13").Value = 679000000.0
$ws.Range("D13").Value = 625000000.0
$ws.Range("E13").Value = 581000000.0
$ws.Range("F13").Value = 611000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = 713000000.0
$ws.Range("C22").Value = 724000000.0
$ws.Range("D22").Value = 761000000.0
$ws.Range("E22").Value = 788000000.0
$ws.Range("F22").Value = 767000000.0

# Row 39: Net Debt
$ws.Range("G39").Value = 4929800000.0

# Row 40: Total Debt
$ws.Range("G40").Value = 5116500000.0
